# Updated cryptos list (price + 1h volume change) to match the latest scrape.
# D-column "price" text looks numeric for most rows (e.g. "1.012"), so a plain
# Range.Value assignment would make Excel silently coerce it to a real number
# cell. Prefixing the literal with an apostrophe forces Excel to keep it as text
# (matching the original inlineStr cells), and re-applying the "Normal" style
# afterwards clears the quote-prefix formatting flag picked up along the way so
# the cell keeps its original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.965.62'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '''1.843.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("D4").Value = '''1.012'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").Value = '''1.012'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '

$ws.Range("D6").Value = '''308.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.51%  '

$ws.Range("D7").Value = '''0.4773'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.93%  '

$ws.Range("D8").Value = '''0.3677'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.96%  '

$ws.Range("D9").Value = '''0.07210'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.36%  '

$ws.Range("D10").Value = '''0.9293'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.59%  '

$ws.Range("D11").Value = '''19.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.85%  '

$ws.Range("D12").Value = '''0.07717'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("D13").Value = '''1.851.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.68%  '

$ws.Range("D14").Value = '''5.427'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.47%  '

$ws.Range("D15").Value = '''6.449'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.83%  '

$ws.Range("D16").Value = '''88.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.41%  '

$ws.Range("E17").Value = '  +0.34%  '

$ws.Range("D18").Value = '''0.000008647'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.50%  '

$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("D20").Value = '''27.015.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.02%  '

$ws.Range("E21").Value = '  +1.17%  '

$ws.Range("E22").Value = '  +0.48%  '

$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("D24").Value = '''1.944'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.77%  '

$ws.Range("D25").Value = '''152.53'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").Value = '''18.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.73%  '

$ws.Range("D27").Value = '''2.011'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.82%  '

$ws.Range("D28").Value = '''114.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.20%  '

$ws.Range("D29").Value = '''4.958'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.29%  '

$ws.Range("D30").Value = '''0.08863'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.12%  '

$ws.Range("E31").Value = '  +4.15%  '

$ws.Range("D32").Value = '''1.176'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.52%  '

$ws.Range("D33").Value = '''0.7400'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.26%  '

$ws.Range("D34").Value = '''4.492'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.39%  '

$ws.Range("D35").Value = '''2.694'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.02%  '

$ws.Range("D36").Value = '''1.108'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.58%  '

$ws.Range("D37").Value = '''0.01959'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.10%  '

$ws.Range("D38").Value = '''0.05248'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.55%  '

$ws.Range("D39").Value = '''2.961'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("D40").Value = '''0.5249'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.22%  '

$ws.Range("D41").Value = '''7.003'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.13%  '

$ws.Range("D42").Value = '''0.1511'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("D43").Value = '''8.280'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.12%  '

$ws.Range("D44").Value = '''10.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.86%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").Value = '''1.013'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.42%  '

$ws.Range("E47").Value = '  +1.38%  '

$ws.Range("D48").Value = '''1.604'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '

$ws.Range("D49").Value = '''65.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.22%  '

$ws.Range("D50").Value = '''0.06069'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.21%  '

$ws.Range("D51").Value = '''0.8885'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.12%  '
